$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Install K8s - 3"): master-node init is now done as root, so the
# explicit "sudo" before "kubeadm init ..." is no longer needed.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item(2)
$tr10 = $shape10.TextFrame.TextRange

# "- Initialize Kubernetes on master node" -> "...as root" (single run)
$titlePara = $tr10.Paragraphs(1, 1)
$titleRange = $tr10.Characters($titlePara.Start, $titlePara.Length)
$titleRange.Text = "- Initialize Kubernetes on master node as root"

# " sudo kubeadm init --pod-network-cidr=10.244.0.0/16" -> drop the "sudo "
$cmdPara = $tr10.Paragraphs(2, 1)
$sudoRange = $tr10.Characters($cmdPara.Start + 1, 5)
if ($sudoRange.Text -eq "sudo ") {
    $sudoRange.Delete()
}

# ---------------------------------------------------------------------------
# Slide 8 ("Install K8s - 1"): switch the apt repo URL to https and keep the
# trailing slash attached to the hyperlinked text instead of as its own run.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$tr8 = $shape8.TextFrame.TextRange

$repoPara = $tr8.Paragraphs(4, 1)
$fullText = $repoPara.Text
$urlText = "http://apt.kubernetes.io"
$idx = $fullText.IndexOf($urlText)
if ($idx -ge 0) {
    $urlRange = $tr8.Characters($repoPara.Start + $idx, $urlText.Length)
    $urlRange.Text = "https://apt.kubernetes.io/"

    $repoPara2 = $tr8.Paragraphs(4, 1)
    $slashPos = $repoPara2.Start + $idx + "https://apt.kubernetes.io/".Length
    $slashRange = $tr8.Characters($slashPos, 1)
    if ($slashRange.Text -eq "/") {
        $slashRange.Text = " "
    }
}
